# Weekly update: a new "Arveja Verde" price-report row for Femacal de La
# Calera is inserted right after the existing block's first 4 rows (i.e.
# as the new row 6), pushing every subsequent row down by one. The sheet's
# used range therefore grows from A1:R84 to A1:R85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 6 - Excel shifts rows 6:84 down to 7:85
# and copies formatting from the row above, matching native Excel
# "Insert Sheet Rows" behaviour.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with this week's reading.
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 'Femacal de La Calera'
$ws.Range("C6").Value = 'Coquimbo'
$ws.Range("D6").Value = 45111
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 100112022
$ws.Range("G6").Value = 'Arveja Verde'
$ws.Range("H6").Value = 'Perfection'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 35
$ws.Range("K6").Value = 24000
$ws.Range("L6").Value = 24000
$ws.Range("M6").Value = 24000
$ws.Range("N6").Value = '$/saco 25 kilos'
$ws.Range("O6").Value = 'Provincia de Limarí'
$ws.Range("P6").Value = 960
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = 'Hortaliza'
